$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-11: columns C (loss) and D (accuracy) get new values ---
$ws.Range("C2").Value = 0.6793167591094971
$ws.Range("D2").Value = 0.5790757536888123
$ws.Range("C3").Value = 0.6305791735649109
$ws.Range("D3").Value = 0.6393302083015442
$ws.Range("C4").Value = 0.5334708094596863
$ws.Range("D4").Value = 0.7329421639442444
$ws.Range("C5").Value = 0.4828438460826874
$ws.Range("D5").Value = 0.7683433294296265
$ws.Range("C6").Value = 0.4596393704414368
$ws.Range("D6").Value = 0.7874612808227539
$ws.Range("C7").Value = 0.4471473097801208
$ws.Range("D7").Value = 0.7941195964813232
$ws.Range("C8").Value = 0.4223135411739349
$ws.Range("D8").Value = 0.8045355677604675
$ws.Range("C9").Value = 0.4200760722160339
$ws.Range("D9").Value = 0.8075680732727051
$ws.Range("C10").Value = 0.4121015965938568
$ws.Range("D10").Value = 0.8127101063728333
$ws.Range("C11").Value = 0.3990113139152527
$ws.Range("D11").Value = 0.8204891681671143

# --- Append new epochs 11-30 as rows 12-31 ---
# Copy row-11 formatting (bold/centered/bordered A-column style) down through row 31 first
$ws.Range("A11").Copy($ws.Range("A12:A31"))

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 11
$ws.Range("C12").Value = 0.3999518156051636
$ws.Range("D12").Value = 0.8215439319610596
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = 0.3939414024353027
$ws.Range("D13").Value = 0.8257630467414856
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 13
$ws.Range("C14").Value = 0.3930386304855347
$ws.Range("D14").Value = 0.824906051158905
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 14
$ws.Range("C15").Value = 0.3911015391349792
$ws.Range("D15").Value = 0.824114978313446
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 15
$ws.Range("C16").Value = 0.3888201713562012
$ws.Range("D16").Value = 0.8231920599937439
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 16
$ws.Range("C17").Value = 0.3787821233272552
$ws.Range("D17").Value = 0.8309051394462585
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 17
$ws.Range("C18").Value = 0.3797362744808197
$ws.Range("D18").Value = 0.8326850533485413
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 18
$ws.Range("C19").Value = 0.3773754239082336
$ws.Range("D19").Value = 0.8341354131698608
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 19
$ws.Range("C20").Value = 0.3776697814464569
$ws.Range("D20").Value = 0.8306414484977722
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 20
$ws.Range("C21").Value = 0.3735400140285492
$ws.Range("D21").Value = 0.8334102630615234
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 21
$ws.Range("C22").Value = 0.3717670738697052
$ws.Range("D22").Value = 0.8351901769638062
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 22
$ws.Range("C23").Value = 0.3700553774833679
$ws.Range("D23").Value = 0.8362450003623962
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = 23
$ws.Range("C24").Value = 0.368677020072937
$ws.Range("D24").Value = 0.839277446269989
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 24
$ws.Range("C25").Value = 0.3648383617401123
$ws.Range("D25").Value = 0.8412551879882812
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 25
$ws.Range("C26").Value = 0.3613282740116119
$ws.Range("D26").Value = 0.8415188789367676
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = 26
$ws.Range("C27").Value = 0.3652969896793365
$ws.Range("D27").Value = 0.8402663469314575
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = 27
$ws.Range("C28").Value = 0.3593234717845917
$ws.Range("D28").Value = 0.8412551879882812
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = 28
$ws.Range("C29").Value = 0.3594219386577606
$ws.Range("D29").Value = 0.8403322696685791
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = 29
$ws.Range("C30").Value = 0.3544694483280182
$ws.Range("D30").Value = 0.844880998134613
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = 30
$ws.Range("C31").Value = 0.3567031621932983
$ws.Range("D31").Value = 0.842837393283844
